$d = $word.ActiveDocument

# Change 1: fill in the "jabatan" placeholder with the literal text "Wali Nagari"
$d.Content.Find.Execute("{jabatan_orang_1}", $true, $false, $false, $false, $false, $true, 1, $false, "Wali Nagari", 2)

# Change 2a: drop the leading "AN " run so the signature line just reads
# "WALI NAGARI LIMO KOTO"
$d.Content.Find.Execute("AN WALI NAGARI LIMO KOTO", $true, $false, $false, $false, $false, $true, 1, $false, "WALI NAGARI LIMO KOTO", 2)

# Change 2b: remove the now-unwanted "Sekretaris " paragraph entirely
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -match "^Sekretaris\s*$") {
        $p.Range.Delete()
        break
    }
}
